# Insert a new weekly record at row 18 (Macroferia Regional de Talca - Esparragos),
# pushing the existing rows 18..102 down to 19..103.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 18; this shifts rows 18-102 down to 19-103
# and Excel auto-extends the sheet dimension to A1:R103.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new weekly data point.
$ws.Range("A18").Value = 5
$ws.Range("B18").Value = "Macroferia Regional de Talca"
$ws.Range("C18").Value = "Maule"
$ws.Range("D18").Value = "09/20/2023"
$ws.Range("E18").Value = 7
$ws.Range("F18").Value = 300000000
$ws.Range("G18").Value = "Espárragos"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 2000
$ws.Range("L18").Value = 2000
$ws.Range("M18").Value = 2000
$ws.Range("N18").Value = "$/kilo"
$ws.Range("O18").Value = "Provincia de Linares"
$ws.Range("P18").Value = 2000
$ws.Range("Q18").Value = 1
$ws.Range("R18").Value = "Hortaliza"
